$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style from H1 (bold header style) and set values
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J1").Value = "IF"

$excel.CutCopyMode = 0

# Data rows 2-74: I and J column numeric values
$data = @(
    @(2,6,7),
    @(3,6,6),
    @(4,6,6),
    @(5,6,6),
    @(6,7,7),
    @(7,8,9),
    @(8,5,5),
    @(9,6,7),
    @(10,8,8),
    @(11,6,7),
    @(12,7,7),
    @(13,7,7),
    @(14,5,6),
    @(15,8,8),
    @(16,6,6),
    @(17,5,6),
    @(18,5,6),
    @(19,10,11),
    @(20,7,9),
    @(21,6,6),
    @(22,7,7),
    @(23,9,9),
    @(24,7,7),
    @(25,10,10),
    @(26,12,12),
    @(27,5,5),
    @(28,6,7),
    @(29,9,9),
    @(30,9,9),
    @(31,9,9),
    @(32,9,9),
    @(33,9,9),
    @(34,9,9),
    @(35,9,9),
    @(36,9,9),
    @(37,9,9),
    @(38,9,9),
    @(39,9,9),
    @(40,9,9),
    @(41,9,9),
    @(42,9,9),
    @(43,9,9),
    @(44,8,9),
    @(45,7,7),
    @(46,9,9),
    @(47,8,8),
    @(48,7,7),
    @(49,8,8),
    @(50,7,7),
    @(51,5,5),
    @(52,9,9),
    @(53,8,8),
    @(54,9,9),
    @(55,6,6),
    @(56,6,7),
    @(57,6,6),
    @(58,6,6),
    @(59,8,8),
    @(60,9,9),
    @(61,8,8),
    @(62,9,9),
    @(63,8,8),
    @(64,8,8),
    @(65,9,10),
    @(66,9,9),
    @(67,9,9),
    @(68,10,10),
    @(69,9,9),
    @(70,8,8),
    @(71,4,4),
    @(72,5,5),
    @(73,4,4),
    @(74,3,3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
